# Populate the "Durable" sheet (sheet3) with the same kind of experiment
# table already present on "Basic_problemoverflow" / "Basic", and nudge the
# saved selections, matching the target revision.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)

$xlCenter = -4108

# ---- Row 1 / Row 2 header cells -----------------------------------------
$headerRow1 = @{
    A = "TLC模型"
    D = "状态图直径"
    E = "状态数"
    F = "不同状态数"
    G = "Queue Size"
    H = "检验时间"
}
foreach ($col in @("A","B","C","D","E","F","G","H")) {
    $cell = $ws3.Range($col + "1")
    if ($headerRow1.ContainsKey($col)) {
        $cell.Value = $headerRow1[$col]
    }
    $cell.HorizontalAlignment = $xlCenter
    $cell.Font.Color = 0
}

$headerRow2 = @{
    A = "Server Num"
    B = "Client Num"
    C = "Key Num"
}
foreach ($col in @("A","B","C","D","E","F","G","H")) {
    $cell = $ws3.Range($col + "2")
    if ($headerRow2.ContainsKey($col)) {
        $cell.Value = $headerRow2[$col]
    }
    $cell.HorizontalAlignment = $xlCenter
    $cell.Font.Color = 0
}

# ---- Data rows 3-18 -------------------------------------------------------
# Columns: A,B,C,D,E,F,G are plain numbers; H is a duration (h:mm:ss).
$rows = @(
    @(3,2,2,15,351846374,100000027,68427817,0.04297453703703704),
    @(3,2,3,15,349869869,100000016,69236467,0.044016203703703703),
    @(3,2,4,14,347023882,100000029,69791592,0.042465277777777775),
    @(3,2,5,14,346545192,100000012,70234362,0.042928240740740746),
    @(3,3,2,13,364402539,100000041,73860653,0.065995370370370371),
    @(3,3,3,13,373571942,100000035,74975922,0.071319444444444449),
    @(3,3,4,13,381596426,100000054,75905896,0.073124999999999996),
    @(3,3,5,12,379711195,100000017,76533973,0.078032407407407411),
    @(5,2,2,13,377009534,100000037,79540261,0.32546296296296295),
    @(5,2,3,$null,$null,$null,$null,$null),
    @(5,2,4,$null,$null,$null,$null,$null),
    @(5,2,5,13,414728822,100000042,79483016,0.41244212962962962),
    @(5,3,2,13,463905335,100000022,80998636,0.89020833333333327),
    @(5,3,3,$null,$null,$null,$null,$null),
    @(5,3,4,$null,$null,$null,$null,$null),
    @(5,3,5,$null,$null,$null,$null,$null)
)

$cols = @("A","B","C","D","E","F","G","H")
$r = 3
foreach ($row in $rows) {
    for ($i = 0; $i -lt 8; $i++) {
        $cell = $ws3.Range($cols[$i] + $r)
        $cell.Value = $row[$i]
        $cell.HorizontalAlignment = $xlCenter
        if ($i -eq 7) {
            $cell.NumberFormat = "h:mm:ss"
        }
        $cell.Font.Color = 0
    }
    $r++
}

# ---- Merged header cells ---------------------------------------------------
$ws3.Range("A1:C1").MergeCells = $true
$ws3.Range("D1:D2").MergeCells = $true
$ws3.Range("E1:E2").MergeCells = $true
$ws3.Range("F1:F2").MergeCells = $true
$ws3.Range("G1:G2").MergeCells = $true
$ws3.Range("H1:H2").MergeCells = $true

# ---- Selections -------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("H10").Select()

$ws3.Activate()
$ws3.Range("J15").Select()
